# summer 24 week 7 updates
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$ws.Range("F4").Value = 10.16

$ws.Range("F5").Value = 10.18
$ws.Range("G5").Value = 9.72

$ws.Range("D6").Value = 9.84
$ws.Range("E6").Value = 9.82
$ws.Range("G6").Value = 10.34
$ws.Range("H6").Value = 10.48

$ws.Range("E7").Value = 10.28
$ws.Range("F7").Value = 9.66
$ws.Range("H7").Value = 9.98

$ws.Range("F8").Value = 9.52
$ws.Range("G8").Value = 10.02
$ws.Range("I8").Value = 8.8

$ws.Range("H9").Value = 11.2
